$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current (before) layout:
# Row1: Name | CAN-ID | Info | Laenge
# Row2: CAN_ID_CONTROL_LOWER_MOTOR | 0xC0 | Kontrolliert das an den unteren Motor ausgegebene PWM-Signal | 1
# Row3: CAN_ID_INFOS_LOWER_CONTROLLER | 0xC1 | Infos des unteren Motors | 6
# Row4: CAN_ID_CONTROL_BACK_MOTOR | 0xC2 | Kontrolliert das an den unteren Motor ausgegebene PWM-Signal | 1
# Row5: CAN_ID_INFOS_BACK_CONTROLLER | 0xC3 | Infos des unteren Motors | 6
# Row6: CAN_ID_CONTROL_SERVOS | 0xD0 | Kontrolliert die fuer die Lenkung zustaendigen Servos | 1
#
# New (after) layout: the two "control" entries for the lower motor and the
# back motor are merged into a single entry (CAN_ID_CONTROL_MOTORS_SERVOS)
# that now also controls the steering servos, keeping CAN-ID 0xC0 but with a
# length of 3. The old CAN_ID_CONTROL_SERVOS entry moves up to take the
# freed row and keeps its own CAN-ID (0xD0), shrinking the table from 6 to 5
# rows.

$ws.Range("A2").Value = "CAN_ID_CONTROL_MOTORS_SERVOS"
$ws.Range("B2").Value = "0xC0"
$ws.Range("C2").Value = "Kontrolliert das an den unteren Motor ausgegebene PWM-Signal"
$ws.Range("D2").Value = 3

$ws.Range("A3").Value = "CAN_ID_INFOS_LOWER_CONTROLLER"
$ws.Range("B3").Value = "0xC1"
$ws.Range("C3").Value = "Infos des unteren Motors"
$ws.Range("D3").Value = 6

$ws.Range("A4").Value = "CAN_ID_INFOS_BACK_CONTROLLER"
$ws.Range("B4").Value = "0xC3"
$ws.Range("C4").Value = "Infos des unteren Motors"
$ws.Range("D4").Value = 6

$ws.Range("A5").Value = "CAN_ID_CONTROL_SERVOS"
$ws.Range("B5").Value = "0xD0"
$ws.Range("C5").Value = "Kontrolliert die für die Lenkung zuständigen Servos"
$ws.Range("D5").Value = 1

# Row 6 no longer exists - clear it out entirely.
$ws.Range("A6:D6").ClearContents()

# Update the selection to match the saved view state.
$ws.Range("B5").Select()

$wb.Save()
